$wb = $excel.ActiveWorkbook

# --- Update status text from "Ready for handoff" to "In Translation" ---
# This shared string is referenced from the Overview sheet (columns E/F, row 2)
# as well as from the zh-cn and de-de sheets (column C, row 2).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the "status" columns ---
# Target stored column width is 13.4101845877511 characters; the closest
# value reachable through the ColumnWidth rounding (6px Maximum Digit Width,
# 5px padding) is produced by requesting a ColumnWidth of 12.5.
$narrowColumnWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $narrowColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $narrowColumnWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $narrowColumnWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $narrowColumnWidth
